$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# TC2 (rows 15-20) and TC3 (rows 23-28) swap their second scenario step.
# Before:
#   TC2 step 2 (row 20): "Chefe Clica para ordenar pelo nome do servidor." / "SYSTEM Visualiza ... ordenado pelo nome do servidor."
#   TC3 step 2 (row 28): "Chefe Indica alguns parâmetros específicos para a busca; ..." / "SYSTEM Exibe uma nova listagem ..."
# After:
#   TC2 step 2 (row 20): "Chefe Indica alguns parâmetros específicos para a busca; ..." / "SYSTEM Exibe uma nova listagem ..."
#   TC3 step 2 (row 28): "Chefe Clica para ordenar pelo nome do servidor." / "SYSTEM Visualiza ... ordenado pelo nome do servidor."

$tc2Steps = $ws.Range("B20").Value2
$tc2Results = $ws.Range("D20").Value2
$tc3Steps = $ws.Range("B28").Value2
$tc3Results = $ws.Range("D28").Value2

$ws.Range("B20").Value2 = $tc3Steps
$ws.Range("D20").Value2 = $tc3Results
$ws.Range("B28").Value2 = $tc2Steps
$ws.Range("D28").Value2 = $tc2Results
